# إضافة حدث جديد في Card18
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# Fill the previously-blank "nan" placeholder cells on row 31
$ws.Range("B31").Value = "nan"
$ws.Range("C31").Value = "nan"
$ws.Range("D31").Value = "nan"
$ws.Range("E31").Value = "nan"
$ws.Range("F31").Value = "nan"
$ws.Range("G31").Value = "nan"
$ws.Range("H31").Value = "nan"
$ws.Range("I31").Value = "nan"
$ws.Range("J31").Value = "nan"
$ws.Range("K31").Value = "nan"
$ws.Range("N31").Value = "nan"
$ws.Range("Q31").Value = "nan"

# Append the new event row (row 32)
# Force text formatting first so numeric-looking / date-looking strings
# ("18", "1/2/2026") are stored as literal text, matching the rest of
# the sheet (every data cell on this sheet is stored as text).
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = "18"
$ws.Range("L32").NumberFormat = "@"
$ws.Range("L32").Value = "1/2/2026"
$ws.Range("M32").Value = "انقطاع سير دوفر 1200"
$ws.Range("O32").Value = "تم تغير سير دوفر 1200"
$ws.Range("P32").Value = "عمر"
